$d = $word.ActiveDocument

# 1. Update the title: "TIMEMASTER" -> "TIMEMASTER PRO" (curly quotes preserved)
$openQuote = [char]0x201C
$closeQuote = [char]0x201D
$d.Content.Find.Execute(
    ($openQuote + "TIMEMASTER" + $closeQuote),
    $true, $false, $false, $false, $false, $true, 1, $false,
    ($openQuote + "TIMEMASTER PRO" + $closeQuote),
    2)

# 2. Replace the author name "Vilcarano De la cruz Frank" with "Sanchez Ramos Giovany Elver"
$d.Content.Find.Execute(
    "Vilcarano De la cruz Frank",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Sanchez Ramos Giovany Elver",
    2)
